$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 173497
$ws.Range("C4").Value = 163861
$ws.Range("C5").Value = 9636
$ws.Range("C6").Value = 739
$ws.Range("C7").Value = 5.55
$ws.Range("C8").Value = 65.97
